$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 10 ("chemical_recycling_gasification") so that it
# shifts the existing rows 10-24 ("fossil_routes" ... "fossil_lock_in") down to 11-25.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new "chemical_recycling_pyrolysis" parameter,
# enabled (TRUE) like the other route toggles.
$ws.Cells.Item(10, 1).Value = "chemical_recycling_pyrolysis"
$ws.Cells.Item(10, 2).Value = $true
